$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 (this shifts rows 19:54 down to 20:55
# and extends the sheet dimension from A1:R54 to A1:R55).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly price record.
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44571
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112021
$ws.Range("G19").Value = "Ají"
$ws.Range("H19").Value = "Americana (o)"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15500
$ws.Range("N19").Value = "$/caja 15 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1033
$ws.Range("Q19").Value = 15
$ws.Range("R19").Value = "Hortaliza"
